$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; this shifts the existing rows 21..119
# down to 22..120 (matching the dimension change from A1:R119 to A1:R120).
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new "Puerro" record.
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 45030
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112005
$ws.Range("G21").Value = "Puerro"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 9000
$ws.Range("M21").Value = 9000
$ws.Range("N21").Value = "$/paquete 20 unidades"
$ws.Range("O21").Value = "Provincia de Chacabuco"
$ws.Range("P21").Value = 450
$ws.Range("Q21").Value = 20
$ws.Range("R21").Value = "Hortaliza"
